# Insert a new weekly price-record row for Níspero (Provincia de Los Andes,
# 2023-11-20) above the current row 9, pushing the existing rows 9-35 down
# to rows 10-36 (the sheet's used range grows from A1:T35 to A1:T36).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 9; Excel shifts rows 9..35 down to 10..36
# and copies formatting (e.g. the date style on column D) from the row above.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new record.
$ws.Cells.Item(9, 1).Value  = 10
$ws.Cells.Item(9, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(9, 3).Value  = "La Araucanía"
$ws.Cells.Item(9, 4).Value  = (Get-Date -Year 2023 -Month 11 -Day 20 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(9, 5).Value  = 9
$ws.Cells.Item(9, 6).Value  = "Fruta"
$ws.Cells.Item(9, 7).Value  = 100104
$ws.Cells.Item(9, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(9, 9).Value  = 100104004
$ws.Cells.Item(9, 10).Value = "Níspero"
$ws.Cells.Item(9, 11).Value = "Californiana(o)"
$ws.Cells.Item(9, 12).Value = "Primera"
$ws.Cells.Item(9, 13).Value = 25
$ws.Cells.Item(9, 14).Value = 30000
$ws.Cells.Item(9, 15).Value = 30000
$ws.Cells.Item(9, 16).Value = 30000
$ws.Cells.Item(9, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(9, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(9, 19).Value = 3000
$ws.Cells.Item(9, 20).Value = 10
